# Updated symbol list on Tue Dec 27 07:24:31 UTC 2022 with GitHub Actions
#
# Applies the cell-level text updates described by the diff to the
# "cryptos" worksheet. All touched columns (B, C, D, E) store plain text
# in the source workbook (t="inlineStr"), so every write below forces a
# Text number format before assigning the value (otherwise numeric-looking
# strings like "243.65" get auto-coerced into real numbers by Excel), and
# then resets the cell style back to "Normal" so no stray number-format
# style is left behind on cells that originally had none.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, [string]$Cell, [string]$Value)
    $rng = $Sheet.Range($Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $Value
    $rng.Style = "Normal"
}

# Row 2 - BNB
Set-TextValue $ws "D2" "243.65"

# Row 3 - OKB
Set-TextValue $ws "D3" "23.05"

# Row 4 - HuobiToken
Set-TextValue $ws "D4" "5.405"

# Row 5 - Cronos
Set-TextValue $ws "D5" "0.05988"

# Row 6
Set-TextValue $ws "D6" "3.422"

# Row 7
Set-TextValue $ws "D7" "6.500"

# Row 8
Set-TextValue $ws "D8" "0.8122"

# Row 10
Set-TextValue $ws "D10" "0.1426"

# Row 11
Set-TextValue $ws "D11" "0.07442"

# Row 12
Set-TextValue $ws "D12" "0.03328"

# Row 13
Set-TextValue $ws "D13" "0.03066"

# Row 14
Set-TextValue $ws "D14" "0.09351"

# Row 16
Set-TextValue $ws "D16" "0.001578"

# Row 17
Set-TextValue $ws "D17" "0.04705"

# Row 18 - One
Set-TextValue $ws "D18" "0.01114"
Set-TextValue $ws "E18" "17OneONEBestin24h"

# Row 19 - TigerCash
Set-TextValue $ws "D19" "0.005928"

# Row 20 - BitKan
Set-TextValue $ws "E20" "19BitKanKAN"

# Row 21 - HotbitToken
Set-TextValue $ws "D21" "0.004881"

# Row 22 - NitroEx
Set-TextValue $ws "D22" "0.00008000"
Set-TextValue $ws "E22" "21NitroExNTXWorstin24h"

# Row 40 - IDEX
Set-TextValue $ws "D40" "0.03949"

# Row 41 - was BKEXToken, now KickToken
Set-TextValue $ws "B41" "KickToken"
Set-TextValue $ws "C41" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws "D41" "0.006370"
Set-TextValue $ws "E41" "40KickTokenKICK"

# Row 42 - CEJI
Set-TextValue $ws "D42" "0.004000"

# Row 43 - was KickToken, now BKEXToken
Set-TextValue $ws "B43" "BKEXToken"
Set-TextValue $ws "C43" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D43" "0.1077"
Set-TextValue $ws "E43" "42BKEXTokenBKK"

# Row 44
Set-TextValue $ws "D44" "0.009090"

# Row 45
Set-TextValue $ws "D45" "0.00005200"

# Row 48
Set-TextValue $ws "D48" "0.002267"
